# Remove the "GitHub: https://github.com/mdebadwar" paragraph entirely
# (the paragraph that sits between the "Location: ..." line and the
# "LinkedIn: ..." line in the cover-letter header table).

$d = $word.ActiveDocument

$target = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -match "^GitHub:\s*https://github\.com/mdebadwar\s*$") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    # Delete the whole paragraph, including its end-of-paragraph mark,
    # so the following paragraph (LinkedIn:) simply moves up - no blank
    # paragraph is left behind.
    $target.Range.Delete()
}
